$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description cells in row 2 where line breaks were normalized to spaces
$ws.Range("B2").Value = "May write out name or use intial with a period. Examples:    Mary or M."
$ws.Range("I2").Value = "Enter one of these values: Home   Work  Mobile"
$ws.Range("K2").Value = "Enter one of these values: Home        Work      Mobile"
$ws.Range("L2").Value = "State is required if applicable in your country. Spell out the state name. Example: North Carolina"

# Update the selected cell in the sheet view
$ws.Range("E6").Select()
